# Updates cryptocurrency price (D) and volume change (E) columns
# per the scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.849.14"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "1.812.86"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4657"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3682"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07362"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").Value = "1.811.77"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008686"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "26.894.82"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.332"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("D24").Value = "2.054.64"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.901"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.176"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.312"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08930"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7654"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.507"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.901"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.086"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01958"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05280"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.940"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.249"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5313"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.60%  "

$ws.Range("E43").Value = "  -1.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1661"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.414"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4922"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("E49").Value = "  +1.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("E51").Value = "  -0.20%  "
